$wb = $excel.ActiveWorkbook

# --- NAND sheet: fill truth-table results for Q = NOT(A AND B) ---
$nand = $wb.Worksheets.Item("NAND")
$nand.Activate() | Out-Null
$nand.Range("D4").Value = 1
$nand.Range("D5").Value = 1
$nand.Range("D6").Value = 1
$nand.Range("D7").Value = 0
$nand.Range("D13").Select() | Out-Null

# --- NOR sheet: fill truth-table results for Q = NOT(A OR B) ---
$nor = $wb.Worksheets.Item("NOR")
$nor.Activate() | Out-Null
$nor.Range("D4").Value = 1
$nor.Range("D5").Value = 0
$nor.Range("D6").Value = 0
$nor.Range("D7").Value = 0
$nor.Range("F5").Select() | Out-Null

# --- XOR sheet: fill truth-table results for Y = (NOT A AND B) OR (A AND NOT B) ---
$xor = $wb.Worksheets.Item("XOR")
$xor.Activate() | Out-Null
$xor.Range("D4").Value = 0
$xor.Range("D5").Value = 1
$xor.Range("D6").Value = 1
$xor.Range("D7").Value = 0
$xor.Range("D23").Select() | Out-Null
